# Update cryptos list values (Price and Volume(1h) columns) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.101.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.318.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.79%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.60%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("E8").Value = "  -0.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.341.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("E10").Value = "  -1.18%  "

$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.737.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.137.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.75%  "

$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.326.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("E20").Value = "  -1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.84%  "

$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.47%  "

$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0727"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("E36").Value = "  -3.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.916"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "39.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "149.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.376"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.58%  "

$ws.Range("E44").Value = "  -1.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "281.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0928"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0501"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.558"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.53%  "

$ws.Range("E50").Value = "  -1.94%  "

$ws.Range("E51").Value = "  -0.98%  "
